# Gallery_Sounder_PFI.xlsx - "Updated TestData for Portugal Market"
#
# 1. Germany sheet: selection changed to A1:D10 (was whole-sheet select)
# 2. A new "Portugal" worksheet is added after "Swiss", cloned from the
#    "Swiss" sheet (same layout/styles/merges), with Portugal-specific
#    data, narrower columns and taller data rows, and becomes the active
#    (selected) sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Germany: just reselect A1:D10 -----------------------------------
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Range("A1:D10").Select()

# --- Build the new "Portugal" sheet off of a copy of "Swiss" ---------
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy([System.Reflection.Missing]::Value, $swiss)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# Market name / JIRA reference for Portugal
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2412"

# Column widths specific to the Portugal sheet
$portugal.Columns.Item(1).ColumnWidth = 24.44140625
$portugal.Columns.Item(2).ColumnWidth = 14.77734375
$portugal.Columns.Item(3).ColumnWidth = 13.44140625
$portugal.Columns.Item(4).ColumnWidth = 13.88671875

# Rows 3-5 grew taller (double default height) on this sheet
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# Make Portugal the active sheet/tab, with B4:B5 selected
$portugal.Activate()
$portugal.Range("B4:B5").Select()
